# Weekly update: insert two new price records (dated 44438) at the top of
# the Alcachofa / Vega Modelo de Temuco table, pushing the existing rows
# (previously at 63..92) down to 65..94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 63 — this shifts rows 63:92 down to 65:94
# and extends the used range to row 94 automatically.
$ws.Rows("63:64").Insert()

# New row 63: Alcachofa, Española, Primera
$ws.Cells.Item(63, 1).Value  = 10
$ws.Cells.Item(63, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(63, 3).Value  = 'La Araucanía'
$ws.Cells.Item(63, 4).Value  = 44438
$ws.Cells.Item(63, 5).Value  = 9
$ws.Cells.Item(63, 6).Value  = 100112013
$ws.Cells.Item(63, 7).Value  = 'Alcachofa'
$ws.Cells.Item(63, 8).Value  = 'Española'
$ws.Cells.Item(63, 9).Value  = 'Primera'
$ws.Cells.Item(63, 10).Value = 200
$ws.Cells.Item(63, 11).Value = 15000
$ws.Cells.Item(63, 12).Value = 15000
$ws.Cells.Item(63, 13).Value = 15000
$ws.Cells.Item(63, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(63, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(63, 16).Value = 500
$ws.Cells.Item(63, 17).Value = 30
$ws.Cells.Item(63, 18).Value = 'Hortaliza'

# New row 64: Alcachofa, Madrigal, Primera
$ws.Cells.Item(64, 1).Value  = 10
$ws.Cells.Item(64, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(64, 3).Value  = 'La Araucanía'
$ws.Cells.Item(64, 4).Value  = 44438
$ws.Cells.Item(64, 5).Value  = 9
$ws.Cells.Item(64, 6).Value  = 100112013
$ws.Cells.Item(64, 7).Value  = 'Alcachofa'
$ws.Cells.Item(64, 8).Value  = 'Madrigal'
$ws.Cells.Item(64, 9).Value  = 'Primera'
$ws.Cells.Item(64, 10).Value = 100
$ws.Cells.Item(64, 11).Value = 13000
$ws.Cells.Item(64, 12).Value = 13000
$ws.Cells.Item(64, 13).Value = 13000
$ws.Cells.Item(64, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(64, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(64, 16).Value = 325
$ws.Cells.Item(64, 17).Value = 40
$ws.Cells.Item(64, 18).Value = 'Hortaliza'
